# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 14:22"

# --- Plain numeric refreshes (country stays the same) ---
# Paises Bajos
$ws.Cells.Item(17,2).Value = 33405
$ws.Cells.Item(17,3).Value = 750
$ws.Cells.Item(17,4).Value = 250
$ws.Cells.Item(17,5).Value = 29404
$ws.Cells.Item(17,6).Value = 1176
$ws.Cells.Item(17,7).Value = 67
$ws.Cells.Item(17,8).Value = 3751

# Suecia
$ws.Cells.Item(24,2).Value = 14777
$ws.Cells.Item(24,3).Value = 392
$ws.Cells.Item(24,4).Value = 550
$ws.Cells.Item(24,5).Value = 12647
$ws.Cells.Item(24,6).Value = 521
$ws.Cells.Item(24,7).Value = 40
$ws.Cells.Item(24,8).Value = 1580

# Dinamarca
$ws.Cells.Item(36,2).Value = 7515
$ws.Cells.Item(36,3).Value = 131
$ws.Cells.Item(36,4).Value = 4141
$ws.Cells.Item(36,5).Value = 3010
$ws.Cells.Item(36,6).Value = 84
$ws.Cells.Item(36,7).Value = 9
$ws.Cells.Item(36,8).Value = 364

# Croacia
$ws.Cells.Item(65,2).Value = 1881
$ws.Cells.Item(65,3).Value = 10
$ws.Cells.Item(65,4).Value = 771
$ws.Cells.Item(65,5).Value = 1063
$ws.Cells.Item(65,6).Value = 18
$ws.Cells.Item(65,7).Value = 0
$ws.Cells.Item(65,8).Value = 47

# Republica de Macedonia
$ws.Cells.Item(78,2).Value = 1225
$ws.Cells.Item(78,3).Value = 18
$ws.Cells.Item(78,4).Value = 200
$ws.Cells.Item(78,5).Value = 971
$ws.Cells.Item(78,6).Value = 9
$ws.Cells.Item(78,7).Value = 3
$ws.Cells.Item(78,8).Value = 54

# Eslovaquia
$ws.Cells.Item(79,2).Value = 1173
$ws.Cells.Item(79,3).Value = 12
$ws.Cells.Item(79,4).Value = 251
$ws.Cells.Item(79,5).Value = 909
$ws.Cells.Item(79,6).Value = 8
$ws.Cells.Item(79,7).Value = 1
$ws.Cells.Item(79,8).Value = 13

# Isla de Man
$ws.Cells.Item(116,2).Value = 299
$ws.Cells.Item(116,3).Value = 1
$ws.Cells.Item(116,4).Value = 200
$ws.Cells.Item(116,5).Value = 93
$ws.Cells.Item(116,6).Value = 11
$ws.Cells.Item(116,7).Value = 0
$ws.Cells.Item(116,8).Value = 6

# Vietnam
$ws.Cells.Item(120,2).Value = 268
$ws.Cells.Item(120,3).Value = 0
$ws.Cells.Item(120,4).Value = 214
$ws.Cells.Item(120,5).Value = 54
$ws.Cells.Item(120,6).Value = 8
$ws.Cells.Item(120,7).Value = 0
$ws.Cells.Item(120,8).Value = 0

# --- Congo inserted right after Martinica (row 129); Guadalupe/Ruanda shift down one row ---
$ws.Cells.Item(130,1).Value = "Congo"
$ws.Cells.Item(130,2).Value = 160
$ws.Cells.Item(130,3).Value = 17
$ws.Cells.Item(130,4).Value = 16
$ws.Cells.Item(130,5).Value = 138
$ws.Cells.Item(130,6).Value = 0
$ws.Cells.Item(130,7).Value = 0
$ws.Cells.Item(130,8).Value = 6

$ws.Cells.Item(131,1).Value = "Guadalupe"
$ws.Cells.Item(131,2).Value = 148
$ws.Cells.Item(131,3).Value = 0
$ws.Cells.Item(131,4).Value = 73
$ws.Cells.Item(131,5).Value = 67
$ws.Cells.Item(131,6).Value = 13
$ws.Cells.Item(131,7).Value = 0
$ws.Cells.Item(131,8).Value = 8

$ws.Cells.Item(132,1).Value = "Ruanda"
$ws.Cells.Item(132,2).Value = 147
$ws.Cells.Item(132,3).Value = 0
$ws.Cells.Item(132,4).Value = 76
$ws.Cells.Item(132,5).Value = 71
$ws.Cells.Item(132,6).Value = 0
$ws.Cells.Item(132,7).Value = 0
$ws.Cells.Item(132,8).Value = 0

# --- Benin inserted right after Uganda (row 158); Maldivas..Macao shift down one row ---
$ws.Cells.Item(159,1).Value = "Benin"
$ws.Cells.Item(159,2).Value = 54
$ws.Cells.Item(159,3).Value = 19
$ws.Cells.Item(159,4).Value = 27
$ws.Cells.Item(159,5).Value = 26
$ws.Cells.Item(159,6).Value = 0
$ws.Cells.Item(159,7).Value = 0
$ws.Cells.Item(159,8).Value = 1

$ws.Cells.Item(160,1).Value = "Maldivas"
$ws.Cells.Item(160,2).Value = 52
$ws.Cells.Item(160,3).Value = 0
$ws.Cells.Item(160,4).Value = 16
$ws.Cells.Item(160,5).Value = 36
$ws.Cells.Item(160,6).Value = 0
$ws.Cells.Item(160,7).Value = 0
$ws.Cells.Item(160,8).Value = 0

$ws.Cells.Item(161,1).Value = "Libia"
$ws.Cells.Item(161,2).Value = 51
$ws.Cells.Item(161,3).Value = 0
$ws.Cells.Item(161,4).Value = 11
$ws.Cells.Item(161,5).Value = 39
$ws.Cells.Item(161,6).Value = 0
$ws.Cells.Item(161,7).Value = 0
$ws.Cells.Item(161,8).Value = 1

$ws.Cells.Item(162,1).Value = "Guinea-Bisau"
$ws.Cells.Item(162,2).Value = 50
$ws.Cells.Item(162,3).Value = 0
$ws.Cells.Item(162,4).Value = 3
$ws.Cells.Item(162,5).Value = 47
$ws.Cells.Item(162,6).Value = 0
$ws.Cells.Item(162,7).Value = 0
$ws.Cells.Item(162,8).Value = 0

$ws.Cells.Item(163,1).Value = "Haiti"
$ws.Cells.Item(163,2).Value = 47
$ws.Cells.Item(163,3).Value = 0
$ws.Cells.Item(163,4).Value = 0
$ws.Cells.Item(163,5).Value = 44
$ws.Cells.Item(163,6).Value = 0
$ws.Cells.Item(163,7).Value = 0
$ws.Cells.Item(163,8).Value = 3

$ws.Cells.Item(164,1).Value = "Macao"
$ws.Cells.Item(164,2).Value = 45
$ws.Cells.Item(164,3).Value = 0
$ws.Cells.Item(164,4).Value = 22
$ws.Cells.Item(164,5).Value = 23
$ws.Cells.Item(164,6).Value = 1
$ws.Cells.Item(164,7).Value = 0
$ws.Cells.Item(164,8).Value = 0

# --- Sierra Leona inserted right after Macao (row 164); Puerto Rico..San Martin(FR) shift down one row ---
$ws.Cells.Item(165,1).Value = "Sierra Leona"
$ws.Cells.Item(165,2).Value = 43
$ws.Cells.Item(165,3).Value = 8
$ws.Cells.Item(165,4).Value = 6
$ws.Cells.Item(165,5).Value = 37
$ws.Cells.Item(165,6).Value = 0
$ws.Cells.Item(165,7).Value = 0
$ws.Cells.Item(165,8).Value = 0

$ws.Cells.Item(166,1).Value = "Puerto Rico"
$ws.Cells.Item(166,2).Value = 39
$ws.Cells.Item(166,3).Value = 0
$ws.Cells.Item(166,4).Value = 1
$ws.Cells.Item(166,5).Value = 36
$ws.Cells.Item(166,6).Value = 0
$ws.Cells.Item(166,7).Value = 0
$ws.Cells.Item(166,8).Value = 2

$ws.Cells.Item(167,1).Value = "Eritrea"
$ws.Cells.Item(167,2).Value = 39
$ws.Cells.Item(167,3).Value = 0
$ws.Cells.Item(167,4).Value = 3
$ws.Cells.Item(167,5).Value = 36
$ws.Cells.Item(167,6).Value = 0
$ws.Cells.Item(167,7).Value = 0
$ws.Cells.Item(167,8).Value = 0

$ws.Cells.Item(168,1).Value = "Mozambique"
$ws.Cells.Item(168,2).Value = 39
$ws.Cells.Item(168,3).Value = 0
$ws.Cells.Item(168,4).Value = 8
$ws.Cells.Item(168,5).Value = 31
$ws.Cells.Item(168,6).Value = 0
$ws.Cells.Item(168,7).Value = 0
$ws.Cells.Item(168,8).Value = 0

$ws.Cells.Item(169,1).Value = "Siria"
$ws.Cells.Item(169,2).Value = 39
$ws.Cells.Item(169,3).Value = 0
$ws.Cells.Item(169,4).Value = 5
$ws.Cells.Item(169,5).Value = 31
$ws.Cells.Item(169,6).Value = 0
$ws.Cells.Item(169,7).Value = 0
$ws.Cells.Item(169,8).Value = 3

$ws.Cells.Item(170,1).Value = "San Martin (Parte Francesa)"
$ws.Cells.Item(170,2).Value = 37
$ws.Cells.Item(170,3).Value = 0
$ws.Cells.Item(170,4).Value = 19
$ws.Cells.Item(170,5).Value = 16
$ws.Cells.Item(170,6).Value = 5
$ws.Cells.Item(170,7).Value = 0
$ws.Cells.Item(170,8).Value = 2
